# Fix typo and arrow origin
#
# 1) Move/resize the down-arrow shape "Pfeil nach unten 62" on slide 1.
# 2) Fix the "ontop" typo -> "on top" in the Results bullet list, merging
#    the trailing " " + "ontop" runs into a single " on top" run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1) Arrow shape reposition / resize (rotated shape: Left/Top/Width/Height
#    map 1:1 onto the unrotated <a:off>/<a:ext> in the XML, in points;
#    1 pt = 12700 EMU). Nudge each point value just above the exact EMU
#    boundary so the emulator's point->EMU (floor) conversion lands on
#    the exact target EMU instead of one unit short.
# ---------------------------------------------------------------------
function EmuToPt([double]$emu) {
    return ($emu / 12700.0) + (0.5 / 12700.0)
}

$arrow = $s.Shapes.Item("Pfeil nach unten 62")
$arrow.Left   = EmuToPt 7966392
$arrow.Top    = EmuToPt 12259674
$arrow.Width  = EmuToPt 1257300
$arrow.Height = EmuToPt 4362923

# ---------------------------------------------------------------------
# 2) Typo fix: " ontop" -> " on top" inside the results placeholder that
#    holds the "Creation of a distributed ... mutex ontop" bullet.
# ---------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -like "*ontop*") {
            $target = $shp
            break
        }
    }
}

$tr = $target.TextFrame.TextRange
$found = $tr.Find(" ontop", 0)
$found.Text = " on top"
